$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new item (STRINGAZOLE) showed up in the shortage report, so insert a
# fresh row right above "TORSERETIC 20MG 30 TABS." (old row 53), pushing
# every row below it down by one.
$ws.Rows(53).Insert()

# Copy the formatting of the row that used to be 53 (now shifted to 54)
# onto the newly inserted blank row so it looks like every other data row.
$ws.Range("A54:Q54").Copy()
$ws.Range("A53:Q53").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new row's data.
$ws.Range("A53").Value = 47
$ws.Range("C53").Value = "STRINGAZOLE 40MG 21 ENTERIC COATED TABLETS"
$ws.Range("H53").Value = "0:2"
$ws.Range("L53").Value = "1"
$ws.Range("N53").Value = "126.00"
$ws.Range("P53").Value = "41.5800"
$ws.Range("Q53").Value = "0:1"

# Renumber the following rows (the "م" column is a simple running index).
For ($r = 54; $r -le 75; $r++) {
  $ws.Cells.Item($r, 1).Value = $r - 6
}

# Update the grand-total cell to include the new row's price.
$ws.Range("P76").Value = $ws.Range("P76").Value() + 41.58

# Refresh the generated timestamp shown in the report footer.
$ws.Range("A77").Value = "Sunday, 10 August, 2025 6:54 PM"
